# Change configuration to use config.json.
# - Row 2 (Test #1): Strategy changes from ScalpEmaRsiAdx_X -> MACD
# - Rows 3-5 added as copies of row 2's layout/format, each a new test case
#   (Test #2/3/4) cycling through the remaining strategies, with the
#   strategy that used to be on row 2 (ScalpEmaRsiAdx_X) now ending up on
#   row 5.
# - Selection moves to I10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy row 2's cell formatting down to rows 3:5 so the new rows pick up the
# same number formats / borders / alignment as the existing test row.
$ws.Range("A2:J2").Copy()
$ws.Range("A3:J5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Test #1 (row 2): swap the strategy used for this backtest.
$ws.Range("J2").Value = "MACD"

# Test #2 (row 3)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Bybit"
$ws.Range("C3").Value = "BTCUSDT"
$ws.Range("D3").Value = 44579
$ws.Range("E3").Value = 44586
$ws.Range("F3").Value = "3m"
$ws.Range("G3").Value = 10000
$ws.Range("H3").Value = 0.3
$ws.Range("I3").Value = 0.2
$ws.Range("J3").Value = "MACD_X"

# Test #3 (row 4)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Bybit"
$ws.Range("C4").Value = "BTCUSDT"
$ws.Range("D4").Value = 44579
$ws.Range("E4").Value = 44586
$ws.Range("F4").Value = "3m"
$ws.Range("G4").Value = 10000
$ws.Range("H4").Value = 0.3
$ws.Range("I4").Value = 0.2
$ws.Range("J4").Value = "ScalpEmaRsiAdx"

# Test #4 (row 5)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Bybit"
$ws.Range("C5").Value = "BTCUSDT"
$ws.Range("D5").Value = 44579
$ws.Range("E5").Value = 44586
$ws.Range("F5").Value = "3m"
$ws.Range("G5").Value = 10000
$ws.Range("H5").Value = 0.3
$ws.Range("I5").Value = 0.2
$ws.Range("J5").Value = "ScalpEmaRsiAdx_X"

# Update the active selection to match the recorded cursor position.
$ws.Range("I10").Select()
